$wb = $excel.ActiveWorkbook

# The "Daten" sheet holds the constraint data table that needs updating.
$ws = $wb.Worksheets.Item("Daten")

# Update max_dist_c (column E) values for rows 2-5:
# constraint 6 was excluded for sb's that are already assigned,
# raising the allowed max distance for those rows.
$ws.Range("E2").Value = 8000
$ws.Range("E3").Value = 8000
$ws.Range("E4").Value = 8000
$ws.Range("E5").Value = 5000

# Move the active cell selection on the sheet from H5 to E5.
$ws.Activate()
$ws.Range("E5").Select()
